# BAJAJ & MAGMA FILES
# Update PAID AMOUNT and PERFORMANCE columns for rows 2-4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 838869
$ws.Range("D2").Value = 7.88

$ws.Range("C3").Value = 293818
$ws.Range("D3").Value = 1.68

$ws.Range("C4").Value = 228878
$ws.Range("D4").Value = 2.16
